# Add ability to support multiple entries per form
# This extends the single "q1_opt1/q1_opt2" question block into three
# repeated blocks (Entry Date + 9 question columns each), and adds two
# more respondent entries' worth of answers in row 3 and row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1): fix up existing q1 headers, add the rest
# ---------------------------------------------------------------------
$ws.Range("M1").Value = "q1_opt3"
$ws.Range("N1").Value = "q2_opt1"
$ws.Range("O1").Value = "q2_opt2"
$ws.Range("P1").Value = "q2_opt3"
$ws.Range("Q1").Value = "q3_opt1"
$ws.Range("R1").Value = "q3_opt2"
$ws.Range("S1").Value = "q3_opt3"

# ---------------------------------------------------------------------
# 2. Row 2 (labels): rename existing two, add the rest
# ---------------------------------------------------------------------
$ws.Range("K2").Value = "Choice 1"
$ws.Range("L2").Value = "Choice 2"
$ws.Range("M2").Value = "Choice 3"
$ws.Range("N2").Value = "Check 1"
$ws.Range("O2").Value = "Check 2"
$ws.Range("P2").Value = "Check 3"
$ws.Range("Q2").Value = "Stressed out?"
$ws.Range("R2").Value = "Relaxed?"
$ws.Range("S2").Value = "In pain?"

# ---------------------------------------------------------------------
# 3. Row 3 (Super Employee's entry): update the entry-date value, fill
#    in the rest of the first question block, then add two more full
#    "Entry Date + 9 answers" blocks for the extra form submissions.
#    The two new "Entry Date" cells (T3, AD3) need to carry the same
#    date number-format as J3 -- copy J3's format over before writing
#    the value so they share the existing date style instead of a
#    newly-minted one.
# ---------------------------------------------------------------------
$ws.Range("J3").Value = 42707.962511499216
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 5

$ws.Range("J3").Copy()
$ws.Range("T3").PasteSpecial(-4122)
$ws.Range("AD3").PasteSpecial(-4122)

$ws.Range("T3").Value = 42707.97142585126
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 1
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 2
$ws.Range("AB3").Value = 3
$ws.Range("AC3").Value = 5

$ws.Range("AD3").Value = 42708.016006341655
$ws.Range("AE3").Value = 1
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 1
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 0
$ws.Range("AK3").Value = 1
$ws.Range("AL3").Value = 1
$ws.Range("AM3").Value = 1

# ---------------------------------------------------------------------
# 4. Row 4 (Employee's entry): update the entry-date + first block
#    answers (K4/L4 changed; M4:S4 new -- Q4:S4 left blank, no q3
#    answers recorded for this submission), then the two extra blocks.
# ---------------------------------------------------------------------
$ws.Range("J4").Value = 42707.96312438869
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 0

# No q3 answers were recorded for this submission -- Q4:S4 stay blank,
# but are still present as (unstyled, valueless) cells, matching the
# row's normal default formatting. PasteSpecial(xlPasteFormats) only
# copies the (default) format, not K2's text, so the cells stay empty.
$ws.Range("K2").Copy()
$ws.Range("Q4:S4").PasteSpecial(-4122)

$ws.Range("J3").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("AD4").PasteSpecial(-4122)

$ws.Range("T4").Value = 42707.97153375884
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 1
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 2
$ws.Range("AB4").Value = 3
$ws.Range("AC4").Value = 4

$ws.Range("AD4").Value = 42708.01611987991
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 1
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 1
$ws.Range("AK4").Value = 3
$ws.Range("AL4").Value = 3
$ws.Range("AM4").Value = 3

# ---------------------------------------------------------------------
# 5. Column widths for the newly-added columns (K..AM). The first
#    question block (K:S) mirrors the width pattern of the original
#    columns; the two new "Entry Date + answers" blocks (T:AC, AD:AM)
#    repeat it.
# ---------------------------------------------------------------------
$narrowCols  = @(11,12,13,14,15,16,18,19)
foreach ($c in $narrowCols) { $ws.Columns.Item($c).ColumnWidth = 9.0 }

$wideCols = @(17)
foreach ($c in $wideCols) { $ws.Columns.Item($c).ColumnWidth = 11.142857142857142 }

$dateCols = @(20,30)
foreach ($c in $dateCols) { $ws.Columns.Item($c).ColumnWidth = 27.714285714285715 }

$tinyCols = @(21,22,23,24,25,26,27,28,29,31,32,33,34,35,36,37,38,39)
foreach ($c in $tinyCols) { $ws.Columns.Item($c).ColumnWidth = 3.4285714285714284 }

Write-Output "Applied multi-entry form support edit"
